$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new experiment row (row 13) with "Exp 17" parameters
$ws.Range("A13").Value = "Exp 17"
$ws.Range("B13").Value = 0.55
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = "Exp 17.png"

# Move the active selection to F14, matching where Excel leaves the
# cursor after the row was completed
$ws.Range("F14").Select()
